$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates for rows 2-49 ---
$ws.Range("D2").Value = "64.967.09"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.949.11"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'569.53"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").Value = "'158.94"
$ws.Range("E6").Value = "  +3.48%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "2.943.35"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'6.65"
$ws.Range("E10").Value = "  -4.81%  "
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "'34.04"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "65.070.04"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "3.440.27"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "'6.92"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "2.948.94"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").Value = "'445.65"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'13.91"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'7.25"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'12.02"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  -5.83%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").Value = "'2.38"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").Value = "'48.92"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  -6.86%  "
$ws.Range("D40").Value = "'43.84"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("D41").Value = "'0.299"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").Value = "'2.83"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'384.71"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("D46").Value = "'0.0350"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "2.719.76"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "'132.62"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  +0.04%  "

# --- Row 50/51 swap: Stellar and ThetaToken exchange places with new data ---
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.15"
$ws.Range("E50").Value = "  +4.64%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  +0.61%  "
